$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'37.572.06"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.82%  "
$ws.Range("D3").Value = "'2.036.76"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.24%  "
$ws.Range("E4").Value = "  -0.69%  "
$ws.Range("D5").Value = "'229.74"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.03%  "
$ws.Range("D6").Value = "'0.614"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.41%  "
$ws.Range("E7").Value = "  -0.09%  "
$ws.Range("D8").Value = "'56.14"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +2.96%  "
$ws.Range("D9").Value = "'0.382"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.05%  "
$ws.Range("D10").Value = "'0.0802"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.81%  "
$ws.Range("E11").Value = "  -0.77%  "
$ws.Range("D12").Value = "'2.338.03"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.06%  "
$ws.Range("D13").Value = "'14.42"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.65%  "
$ws.Range("D14").Value = "'20.28"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.10%  "
$ws.Range("D15").Value = "'5.23"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.32%  "
$ws.Range("D16").Value = "'0.743"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.94%  "
$ws.Range("D17").Value = "'2.028.69"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.78%  "
$ws.Range("D18").Value = "'37.481.76"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.78%  "
$ws.Range("D19").Value = "'6.23"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.56%  "
$ws.Range("D20").Value = "'69.06"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.60%  "
$ws.Range("D21").Value = "'0.0₃0826"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.20%  "
$ws.Range("D22").Value = "'223.54"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.79%  "
$ws.Range("E23").Value = "  +0.06%  "
$ws.Range("E24").Value = "  +1.63%  "
$ws.Range("E25").Value = "  +4.03%  "
$ws.Range("D26").Value = "'165.08"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.05%  "
$ws.Range("D27").Value = "'9.15"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.20%  "
$ws.Range("D28").Value = "'0.132"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +5.71%  "
$ws.Range("D29").Value = "'18.83"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.12%  "
$ws.Range("D30").Value = "'1.33"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.26%  "
$ws.Range("D31").Value = "'0.118"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.17%  "
$ws.Range("D32").Value = "'4.49"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.10%  "
$ws.Range("D33").Value = "'0.0606"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.20%  "
$ws.Range("D34").Value = "'4.48"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.90%  "
$ws.Range("E35").Value = "  +8.21%  "
$ws.Range("E36").Value = "  -0.08%  "
$ws.Range("D37").Value = "'5.78"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +9.56%  "
$ws.Range("D38").Value = "'3.25"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +4.46%  "
$ws.Range("E39").Value = "  -0.29%  "
$ws.Range("D40").Value = "'1.477.09"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.14%  "
$ws.Range("D41").Value = "'0.0215"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.72%  "
$ws.Range("B42").Value = "Cronos"
$ws.Range("C42").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D42").Value = "'0.0937"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.39%  "
$ws.Range("B43").Value = "Aave"
$ws.Range("C43").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D43").Value = "'95.32"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.49%  "
$ws.Range("B44").Value = "HuobiToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D44").Value = "'2.84"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +3.19%  "
$ws.Range("B45").Value = "InjectiveProtocol"
$ws.Range("C45").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D45").Value = "'16.50"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.42%  "
$ws.Range("B46").Value = "FTXToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D46").Value = "'4.23"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +18.24%  "
$ws.Range("E47").Value = "  -2.14%  "
$ws.Range("E48").Value = "  +1.14%  "
$ws.Range("D49").Value = "'7.13"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.30%  "
$ws.Range("D50").Value = "'2.95"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.24%  "
$ws.Range("D51").Value = "'2.220.75"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.65%  "
